# Auto-generated Excel COM-interop script to apply odds updates
# for Jogos_do_Dia_Betfair_Back_Lay_2026-01-05.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AA2").Value = 29
$ws.Range("AB2").Value = 16
$ws.Range("AI2").Value = 30
$ws.Range("AL2").Value = 40
$ws.Range("AM2").Value = 70
$ws.Range("AN2").Value = 27
$ws.Range("AO2").Value = 14
$ws.Range("F2").Value = 3.35
$ws.Range("G2").Value = 3.45
$ws.Range("H2").Value = 2.22
$ws.Range("I2").Value = 2.26
$ws.Range("J2").Value = 3.9
$ws.Range("K2").Value = 3.95
$ws.Range("N2").Value = 4.7
$ws.Range("O2").Value = 1.25
$ws.Range("P2").Value = 2.26
$ws.Range("Q2").Value = 1.76
$ws.Range("R2").Value = 1.48
$ws.Range("S2").Value = 2.94
$ws.Range("T2").Value = 1.66
$ws.Range("U2").Value = 2.42
$ws.Range("V2").Value = 1.8
$ws.Range("X2").Value = 21
$ws.Range("Y2").Value = 12.5
$ws.Range("Z2").Value = 15
$ws.Range("AD3").Value = 14.5
$ws.Range("AG3").Value = 17
$ws.Range("I3").Value = 2.58
$ws.Range("K3").Value = 3.95
$ws.Range("L3").Value = 1.31
$ws.Range("N3").Value = 3.55
$ws.Range("P3").Value = 1.94
$ws.Range("Q3").Value = 1.87
$ws.Range("S3").Value = 3
$ws.Range("V3").Value = 1.64
$ws.Range("Y3").Value = 990
$ws.Range("T4").Value = 1.05
$ws.Range("AH5").Value = 48
$ws.Range("AJ5").Value = 10
$ws.Range("AL5").Value = 70
$ws.Range("AN5").Value = 9.4
$ws.Range("F5").Value = 1.35
$ws.Range("G5").Value = 1.37
$ws.Range("H5").Value = 13.5
$ws.Range("I5").Value = 16.5
$ws.Range("K5").Value = 5.4
$ws.Range("L5").Value = 1.4
$ws.Range("N5").Value = 3.3
$ws.Range("O5").Value = 1.35
$ws.Range("Q5").Value = 2.04
$ws.Range("S5").Value = 3.75
$ws.Range("T5").Value = 2.58
$ws.Range("U5").Value = 1.57
$ws.Range("Y5").Value = 36
$ws.Range("AI6").Value = 250
$ws.Range("AM6").Value = 580
$ws.Range("F6").Value = 1.94
$ws.Range("G6").Value = 2
$ws.Range("K6").Value = 3.85
$ws.Range("Q6").Value = 1.81
$ws.Range("T6").Value = 1.72
$ws.Range("U6").Value = 2.18
$ws.Range("W6").Value = 2
$ws.Range("AE7").Value = 260
$ws.Range("AJ7").Value = 14
$ws.Range("F7").Value = 1.43
$ws.Range("K7").Value = 5.5
$ws.Range("L7").Value = 1.24
$ws.Range("U7").Value = 2.18
$ws.Range("X7").Value = 34
$ws.Range("Y7").Value = 980
$ws.Range("Z7").Value = 190
$ws.Range("AB8").Value = 5.5
$ws.Range("AD8").Value = 110
$ws.Range("AH8").Value = 95
$ws.Range("G8").Value = 1.73
$ws.Range("I8").Value = 8.6
$ws.Range("T8").Value = 2.54
$ws.Range("X8").Value = 10.5
$ws.Range("Y8").Value = 20
$ws.Range("AC9").Value = 10
$ws.Range("AD9").Value = 60
$ws.Range("AF9").Value = 18
$ws.Range("AG9").Value = 21
$ws.Range("AK9").Value = 65
$ws.Range("AL9").Value = 150
$ws.Range("AN9").Value = 55
$ws.Range("L9").Value = 1.45
$ws.Range("Z9").Value = 120
$ws.Range("AB10").Value = 15
$ws.Range("AD10").Value = 48
$ws.Range("AF10").Value = 22
$ws.Range("AH10").Value = 60
$ws.Range("AJ10").Value = 900
$ws.Range("AK10").Value = 44
$ws.Range("AN10").Value = 29
$ws.Range("F10").Value = 1.72
$ws.Range("G10").Value = 1.81
$ws.Range("H10").Value = 5.3
$ws.Range("J10").Value = 3.75
$ws.Range("K10").Value = 3.95
$ws.Range("P10").Value = 1.88
$ws.Range("S10").Value = 3.55
$ws.Range("T10").Value = 1.88
$ws.Range("U10").Value = 1.92
$ws.Range("V10").Value = 1.2
$ws.Range("W10").Value = 2.24
$ws.Range("X10").Value = 980
$ws.Range("Y10").Value = 38
$ws.Range("AA11").Value = 190
$ws.Range("AC11").Value = 8.199999999999999
$ws.Range("AE11").Value = 120
$ws.Range("AG11").Value = 11.5
$ws.Range("AK11").Value = 46
$ws.Range("AL11").Value = 160
$ws.Range("G11").Value = 2.32
$ws.Range("J11").Value = 3.45
$ws.Range("T11").Value = 1.78
$ws.Range("W11").Value = 1.75
$ws.Range("Z11").Value = 48
$ws.Range("AM12").Value = 500
$ws.Range("F12").Value = 2.34
$ws.Range("G12").Value = 2.46
$ws.Range("J12").Value = 3.4
$ws.Range("K12").Value = 3.45
$ws.Range("L12").Value = 1.48
$ws.Range("P12").Value = 1.76
$ws.Range("Q12").Value = 2.12
$ws.Range("S12").Value = 4.3
$ws.Range("U12").Value = 2.06
$ws.Range("W12").Value = 1.69
$ws.Range("AB13").Value = 8.6
$ws.Range("AC13").Value = 14
$ws.Range("AD13").Value = 46
$ws.Range("AF13").Value = 23
$ws.Range("AG13").Value = 20
$ws.Range("AH13").Value = 990
$ws.Range("AJ13").Value = 900
$ws.Range("AK13").Value = 75
$ws.Range("AL13").Value = 190
$ws.Range("AN13").Value = 55
$ws.Range("G13").Value = 2.04
$ws.Range("J13").Value = 3.5
$ws.Range("K13").Value = 3.7
$ws.Range("P13").Value = 1.79
$ws.Range("T13").Value = 1.86
$ws.Range("W13").Value = 1.96
$ws.Range("X13").Value = 24
$ws.Range("Y13").Value = 30
